$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loads")
$ws.Select()

# Shift existing columns (v_nom_pu, p_nom_mw, q_nom_mvar, bus_idx) two
# places to the right, working from the rightmost column first so that
# values are not clobbered before they are copied.
$ws.Range("G1").Value2 = $ws.Range("E1").Value2
$ws.Range("G2").Value2 = $ws.Range("E2").Value2
$ws.Range("F1").Value2 = $ws.Range("D1").Value2
$ws.Range("F2").Value2 = $ws.Range("D2").Value2
$ws.Range("E1").Value2 = $ws.Range("C1").Value2
$ws.Range("E2").Value2 = $ws.Range("C2").Value2
$ws.Range("D1").Value2 = $ws.Range("B1").Value2
$ws.Range("D2").Value2 = $ws.Range("B2").Value2

# Fill in the two newly introduced leading columns.
$ws.Range("B1").Value2 = "v_nom_kv"
$ws.Range("C1").Value2 = "s_base_mva"
$ws.Range("B2").Value2 = 22
$ws.Range("C2").Value2 = 100

# Append two new trailing columns for shunt admittance parameters.
$ws.Range("H1").Value2 = "g_shunt_pu"
$ws.Range("I1").Value2 = "b_shunt_pu"
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0

$ws.Range("J2").Select()

$trafos = $wb.Worksheets.Item("trafos")
$trafos.Select()
$trafos.Range("I7").Select()
